$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "parisk"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = "CRT"
$ws.Range("E4").Value = "RES"
$ws.Range("F4").Value = "1a2deef4-16ae-43c8-afd3-8fd2e076505e"
$ws.Range("G4").Value = "rJr4kfWCb_annotated.xlsx"
$ws.Range("H4").Value = "Overall, the paper does not provide any insight beyond: i tried this, i tried that and this works better than that; a strong reject."
